$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update summary header cells
$ws.Range("E11").Value = 2228050
$ws.Range("C13").Value = 5

# Rebuild the worker / period detail table (rows 16-52)
$data = @(
    ,@("CC","1192774235","LISETH CAROLINA NAVARRO ABRIL","1607",32000,1630298)
    ,@("CC","1192774235","LISETH CAROLINA NAVARRO ABRIL","1608",32000,1630298)
    ,@("CC","1192774235","LISETH CAROLINA NAVARRO ABRIL","1609",32000,1630298)
    ,@("CC","1192774235","LISETH CAROLINA NAVARRO ABRIL","1610",32000,1630298)
    ,@("CC","1192774235","LISETH CAROLINA NAVARRO ABRIL","1611",32000,1630298)
    ,@("CC","1192774235","LISETH CAROLINA NAVARRO ABRIL","1612",32000,1630298)
    ,@("CC","1107047216","HENRY GUSTAVO DUQUE FERNANDEZ","2212",20000,2500000)
    ,@("CC","1107047216","HENRY GUSTAVO DUQUE FERNANDEZ","2301",100000,2500000)
    ,@("CC","1107047216","HENRY GUSTAVO DUQUE FERNANDEZ","2302",90000,2500000)
    ,@("CC","1129534404","REGINA MAR A SANTANA ESCORCIA","2311",70000,1750000)
    ,@("CC","1129534404","REGINA MAR A SANTANA ESCORCIA","2312",70000,1750000)
    ,@("CC","1129534404","REGINA MAR A SANTANA ESCORCIA","2401",70000,1750000)
    ,@("CC","1129534404","REGINA MAR A SANTANA ESCORCIA","2402",70000,1750000)
    ,@("CC","1129534404","REGINA MAR A SANTANA ESCORCIA","2403",70000,1750000)
    ,@("CC","1235039795","ALCIDES JUNIOR BELLO FIORILLO","2403",55691,1437451)
    ,@("CC","1129534404","REGINA MAR A SANTANA ESCORCIA","2404",70000,1750000)
    ,@("CC","1235039795","ALCIDES JUNIOR BELLO FIORILLO","2404",55691,1437451)
    ,@("CC","1129534404","REGINA MAR A SANTANA ESCORCIA","2405",70000,1750000)
    ,@("CC","1235039795","ALCIDES JUNIOR BELLO FIORILLO","2405",55691,1437451)
    ,@("CC","1002323564","ENYER LUIS FONTALVO YEPES","2406",13997,1312261)
    ,@("CC","1129534404","REGINA MAR A SANTANA ESCORCIA","2406",70000,1750000)
    ,@("CC","1129534404","REGINA MAR A SANTANA ESCORCIA","2407",70000,1750000)
    ,@("CC","1129534404","REGINA MAR A SANTANA ESCORCIA","2408",70000,1750000)
    ,@("CC","1235039795","ALCIDES JUNIOR BELLO FIORILLO","2408",52490,1437451)
    ,@("CC","1129534404","REGINA MAR A SANTANA ESCORCIA","2409",70000,1750000)
    ,@("CC","1235039795","ALCIDES JUNIOR BELLO FIORILLO","2409",52490,1437451)
    ,@("CC","1129534404","REGINA MAR A SANTANA ESCORCIA","2410",70000,1750000)
    ,@("CC","1129534404","REGINA MAR A SANTANA ESCORCIA","2411",70000,1750000)
    ,@("CC","1129534404","REGINA MAR A SANTANA ESCORCIA","2412",70000,1750000)
    ,@("CC","1129534404","REGINA MAR A SANTANA ESCORCIA","2501",70000,1750000)
    ,@("CC","1129534404","REGINA MAR A SANTANA ESCORCIA","2502",70000,1750000)
    ,@("CC","1129534404","REGINA MAR A SANTANA ESCORCIA","2503",70000,1750000)
    ,@("CC","1129534404","REGINA MAR A SANTANA ESCORCIA","2504",70000,1750000)
    ,@("CC","1129534404","REGINA MAR A SANTANA ESCORCIA","2505",70000,1750000)
    ,@("CC","1129534404","REGINA MAR A SANTANA ESCORCIA","2506",70000,1750000)
    ,@("CC","1129534404","REGINA MAR A SANTANA ESCORCIA","2507",70000,1750000)
    ,@("CC","1129534404","REGINA MAR A SANTANA ESCORCIA","2508",70000,1750000)
)

$r = 16
foreach ($row in $data) {
    $ws.Cells.Item($r, 2).Value = $row[0]
    $ws.Cells.Item($r, 3).Value = $row[1]
    $ws.Cells.Item($r, 4).Value = $row[2]
    $ws.Cells.Item($r, 5).Value = $row[3]
    $ws.Cells.Item($r, 6).Value = $row[4]
    $ws.Cells.Item($r, 7).Value = $row[5]
    $r = $r + 1
}

Write-Output "Done updating rows 16-52"
